$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.027123
$ws.Range("H2").Value = 0.081369
$ws.Range("I2").Value = 0.0960827240265261
$ws.Range("J2").Value = 0.09608272402652611
$ws.Range("M2").Value = 9.841031333333333
$ws.Range("N2").Value = 29.523094
$ws.Range("O2").Value = 0.1083017349730097
$ws.Range("P2").Value = 0.1125970533891552
$ws.Range("Q2").Value = 0.266918292854
$ws.Range("R2").Value = 2.402264635686
$ws.Range("S2").Value = 0.01040592571300566
$ws.Range("T2").Value = 0.01081863160699022
$ws.Range("G3").Value = 0.027123
$ws.Range("H3").Value = 0.081369
$ws.Range("I3").Value = 0.0960827240265261
$ws.Range("J3").Value = 0.09608272402652611
$ws.Range("O3").Value = 0.3504595127507141
$ws.Range("P3").Value = 0.3643589687437936
$ws.Range("Q3").Value = 0.8637355152359998
$ws.Range("R3").Value = 7.773619637123999
$ws.Range("S3").Value = 0.03367310464609766
$ws.Range("T3").Value = 0.03500860224039958
$ws.Range("G4").Value = 0.027123
$ws.Range("H4").Value = 0.081369
$ws.Range("I4").Value = 0.0960827240265261
$ws.Range("J4").Value = 0.09608272402652611
$ws.Range("M4").Value = 16.16670066666667
$ws.Range("N4").Value = 48.500102
$ws.Range("O4").Value = 0.1779164877830196
$ws.Range("P4").Value = 0.1849727733236046
$ws.Range("Q4").Value = 0.438489422182
$ws.Range("R4").Value = 3.946404799638
$ws.Range("S4").Value = 0.01709470079542467
$ws.Range("T4").Value = 0.01777268793167307
$ws.Range("G5").Value = 0.027123
$ws.Range("H5").Value = 0.081369
$ws.Range("I5").Value = 0.0960827240265261
$ws.Range("J5").Value = 0.09608272402652611
$ws.Range("M5").Value = 10.399077
$ws.Range("N5").Value = 20.798154
$ws.Range("O5").Value = 0.1144430947397913
$ws.Range("P5").Value = 0.07932132236322763
$ws.Range("Q5").Value = 0.282054165471
$ws.Range("R5").Value = 1.692324992826
$ws.Range("S5").Value = 0.01099600428862494
$ws.Range("T5").Value = 0.007621408726045114
$ws.Range("G6").Value = 0.027123
$ws.Range("H6").Value = 0.081369
$ws.Range("I6").Value = 0.0960827240265261
$ws.Range("J6").Value = 0.09608272402652611
$ws.Range("M6").Value = 22.614852
$ws.Range("N6").Value = 67.84455600000001
$ws.Range("O6").Value = 0.2488791697534654
$ws.Range("P6").Value = 0.258749882180219
$ws.Range("Q6").Value = 0.613382630796
$ws.Range("R6").Value = 5.520443677164001
$ws.Range("S6").Value = 0.02391298858337316
$ws.Range("T6").Value = 0.02486139352141813
$ws.Range("I7").Value = 0.9039172759734738
$ws.Range("J7").Value = 0.9039172759734738
$ws.Range("M7").Value = 9.841031333333333
$ws.Range("N7").Value = 29.523094
$ws.Range("O7").Value = 0.1083017349730097
$ws.Range("P7").Value = 0.1125970533891552
$ws.Range("Q7").Value = 2.51108676017
$ws.Range("R7").Value = 22.59978084153
$ws.Range("S7").Value = 0.09789580926000403
$ws.Range("T7").Value = 0.101778421782165
$ws.Range("I8").Value = 0.9039172759734738
$ws.Range("J8").Value = 0.9039172759734738
$ws.Range("O8").Value = 0.3504595127507141
$ws.Range("P8").Value = 0.3643589687437936
$ws.Range("S8").Value = 0.3167864081046164
$ws.Range("T8").Value = 0.3293503665033941
$ws.Range("I9").Value = 0.9039172759734738
$ws.Range("J9").Value = 0.9039172759734738
$ws.Range("M9").Value = 16.16670066666667
$ws.Range("N9").Value = 48.500102
$ws.Range("O9").Value = 0.1779164877830196
$ws.Range("P9").Value = 0.1849727733236046
$ws.Range("Q9").Value = 4.12517617561
$ws.Range("R9").Value = 37.12658558048999
$ws.Range("S9").Value = 0.1608217869875949
$ws.Range("T9").Value = 0.1672000853919315
$ws.Range("I10").Value = 0.9039172759734738
$ws.Range("J10").Value = 0.9039172759734738
$ws.Range("M10").Value = 10.399077
$ws.Range("N10").Value = 20.798154
$ws.Range("O10").Value = 0.1144430947397913
$ws.Range("P10").Value = 0.07932132236322763
$ws.Range("Q10").Value = 2.653480482705
$ws.Range("R10").Value = 15.92088289623
$ws.Range("S10").Value = 0.1034470904511663
$ws.Range("T10").Value = 0.07169991363718251
$ws.Range("I11").Value = 0.9039172759734738
$ws.Range("J11").Value = 0.9039172759734738
$ws.Range("M11").Value = 22.614852
$ws.Range("N11").Value = 67.84455600000001
$ws.Range("O11").Value = 0.2488791697534654
$ws.Range("P11").Value = 0.258749882180219
$ws.Range("Q11").Value = 5.77051871058
$ws.Range("R11").Value = 51.93466839522
$ws.Range("S11").Value = 0.2249661811700922
$ws.Range("T11").Value = 0.2338884886588009
